$wb = $excel.ActiveWorkbook

# Row 6 (ALC) - hunk 0
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 3656.7273
$ws.Range("I6").Value = 332.7143
$ws.Range("K6").Value = 998.1428999999999
$ws.Range("M6").Value = -886.1428999999999

# Row 17 (ALC) - hunk 1
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 23398.334
$ws.Range("J17").Value = 23398.334
$ws.Range("L17").Value = 70195.00199999999
$ws.Range("N17").Value = -70531.00199999999

# Row 19 (ALC) - hunk 2
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 3291
$ws.Range("I19").Value = 3854.4
$ws.Range("K19").Value = 3854.4
$ws.Range("M19").Value = -3679.4

# Row 32 (ALC) - hunk 3
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 5099
$ws.Range("I32").Value = 5333.3335
$ws.Range("J32").Value = 4747.5
$ws.Range("K32").Value = 5333.3335
$ws.Range("L32").Value = 4747.5
$ws.Range("M32").Value = -5007.3335
$ws.Range("N32").Value = -5399.5

# Row 97 (ALC) - hunk 4
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H97").Value = 1705.25
$ws.Range("J97").Value = 1801.5
$ws.Range("L97").Value = 5404.5
$ws.Range("N97").Value = -6396.5

# Row 98 (ALC) - hunk 5
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 13798.4
$ws.Range("I98").Value = 12997.667
$ws.Range("K98").Value = 12997.667
$ws.Range("M98").Value = -11499.667

# Row 104 (ALC) - hunk 6
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H104").Value = 787.5
$ws.Range("I104").Value = 787.5
$ws.Range("J104").Value = 0
$ws.Range("K104").Value = 2362.5
$ws.Range("L104").Value = 0
$ws.Range("M104").Value = -615.5
$ws.Range("N104").Value = ""

# Row 111 (ALC) - hunk 7
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 2099.125
$ws.Range("I111").Value = 2041.8572
$ws.Range("K111").Value = 6125.571599999999
$ws.Range("M111").Value = -3058.571599999999

# Row 122 (ALC) - hunk 8
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 13798.4
$ws.Range("I122").Value = 12997.667
$ws.Range("K122").Value = 38993.001
$ws.Range("M122").Value = -36543.001

# Row 125 (ALC) - hunk 9
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 1268.0869
$ws.Range("I125").Value = 1022
$ws.Range("K125").Value = 9198
$ws.Range("M125").Value = -6738

# Row 131 (ALC) - hunk 10
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H131").Value = 818
$ws.Range("I131").Value = 845.6316
$ws.Range("J131").Value = 293
$ws.Range("K131").Value = 2536.8948
$ws.Range("L131").Value = 879
$ws.Range("M131").Value = 2503.1052
$ws.Range("N131").Value = -10959

# Row 132 (ALC) - hunk 11
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 3367.1162
$ws.Range("I132").Value = 1367.303
$ws.Range("K132").Value = 4101.909000000001
$ws.Range("M132").Value = -1571.909000000001

# Row 137 (ALC) - hunk 12
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 49726.57
$ws.Range("I137").Value = 2181.0557
$ws.Range("J137").Value = 334999.66
$ws.Range("K137").Value = 6543.1671
$ws.Range("L137").Value = 1004998.98
$ws.Range("M137").Value = -3993.1671
$ws.Range("N137").Value = -1010098.98

# Row 34 (ARM) - hunk 13
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H34").Value = 20163.334
$ws.Range("I34").Value = 10000
$ws.Range("J34").Value = 25245
$ws.Range("K34").Value = 10000
$ws.Range("L34").Value = 25245
$ws.Range("M34").Value = -9729
$ws.Range("N34").Value = -25787

# Row 37 (ARM) - hunk 14
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 25005000
$ws.Range("I37").Value = 50000000
$ws.Range("J37").Value = 10000
$ws.Range("K37").Value = 50000000
$ws.Range("L37").Value = 10000
$ws.Range("M37").Value = -49999727
$ws.Range("N37").Value = -10546

# Row 40 (ARM) - hunk 15
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H40").Value = 27249.5
$ws.Range("J40").Value = 27999.334
$ws.Range("L40").Value = 27999.334
$ws.Range("N40").Value = -28351.334

# Row 45 (ARM) - hunk 16
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 14260.148
$ws.Range("I45").Value = 13270.158
$ws.Range("K45").Value = 13270.158
$ws.Range("M45").Value = -12893.158

# Row 61 (ARM) - hunk 17
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4748.1113
$ws.Range("I61").Value = 4596.8
$ws.Range("J61").Value = 4937.25
$ws.Range("K61").Value = 4596.8
$ws.Range("L61").Value = 4937.25
$ws.Range("M61").Value = -4384.8
$ws.Range("N61").Value = -5361.25

# Row 74 (ARM) - hunk 18
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1134.8889
$ws.Range("I74").Value = 813.7778
$ws.Range("K74").Value = 813.7778
$ws.Range("M74").Value = 60.22220000000004

# Row 77 (ARM) - hunk 19
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1134.8889
$ws.Range("I77").Value = 813.7778
$ws.Range("K77").Value = 4068.889
$ws.Range("M77").Value = 299.1110000000003

# Row 110 (ARM) - hunk 20
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 995.25
$ws.Range("I110").Value = 993.6667
$ws.Range("K110").Value = 993.6667
$ws.Range("M110").Value = 1051.3333

# Row 124 (ARM) - hunk 21
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H124").Value = 79309
$ws.Range("J124").Value = 79309
$ws.Range("L124").Value = 79309
$ws.Range("N124").Value = -89129

# Row 125 (ARM) - hunk 22
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").Value = ""

# Row 131 (ARM) - hunk 23
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H131").Value = 27650
$ws.Range("J131").Value = 25180
$ws.Range("L131").Value = 25180
$ws.Range("N131").Value = -35260

# Row 132 (ARM) - hunk 24
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1845.9025
$ws.Range("I132").Value = 1391.1154
$ws.Range("J132").Value = 2634.2
$ws.Range("K132").Value = 4173.3462
$ws.Range("L132").Value = 7902.599999999999
$ws.Range("M132").Value = -1643.3462
$ws.Range("N132").Value = -12962.6

# Row 136 (ARM) - hunk 25
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 4748.1113
$ws.Range("I136").Value = 4596.8
$ws.Range("J136").Value = 4937.25
$ws.Range("K136").Value = 13790.4
$ws.Range("L136").Value = 14811.75
$ws.Range("M136").Value = -11240.4
$ws.Range("N136").Value = -19911.75

# Row 141 (ARM) - hunk 26
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H141").Value = 79461.5
$ws.Range("J141").Value = 79461.5
$ws.Range("L141").Value = 79461.5
$ws.Range("N141").Value = -89821.5

# Row 6 (BSM) - hunk 27
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H6").Value = 50000
$ws.Range("J6").Value = 50000
$ws.Range("L6").Value = 50000
$ws.Range("N6").Value = -50226

# Row 35 (BSM) - hunk 28
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 30000
$ws.Range("J35").Value = 30000
$ws.Range("L35").Value = 30000
$ws.Range("N35").Value = -30620

# Row 86 (BSM) - hunk 29
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1462.8158
$ws.Range("I86").Value = 1232.6129
$ws.Range("K86").Value = 1232.6129
$ws.Range("M86").Value = -109.6129000000001

# Row 89 (BSM) - hunk 30
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 1462.8158
$ws.Range("I89").Value = 1232.6129
$ws.Range("K89").Value = 6163.0645
$ws.Range("M89").Value = -547.0645000000004

# Row 105 (BSM) - hunk 31
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 4002271.2
$ws.Range("I105").Value = 4547580.5
$ws.Range("K105").Value = 4547580.5
$ws.Range("M105").Value = -4545833.5

# Row 116 (BSM) - hunk 32
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").Value = ""

# Row 134 (BSM) - hunk 33
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1458.8
$ws.Range("I134").Value = 1458.8
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 4376.4
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -1841.4
$ws.Range("N134").Value = ""

# Row 12 (CRP) - hunk 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 570
$ws.Range("J12").Value = 200
$ws.Range("L12").Value = 200
$ws.Range("N12").Value = -540

# Row 20 (CRP) - hunk 35
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 99332
$ws.Range("J20").Value = 99332
$ws.Range("L20").Value = 99332
$ws.Range("N20").Value = -99804

# Row 30 (CRP) - hunk 36
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H30").Value = 99332
$ws.Range("J30").Value = 99332
$ws.Range("L30").Value = 99332
$ws.Range("N30").Value = -99514

# Row 31 (CRP) - hunk 37
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3852.6875
$ws.Range("I31").Value = 1437
$ws.Range("K31").Value = 1437
$ws.Range("M31").Value = -1142

# Row 34 (CRP) - hunk 38
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3852.6875
$ws.Range("I34").Value = 1437
$ws.Range("K34").Value = 1437
$ws.Range("M34").Value = -1235

# Row 58 (CRP) - hunk 39
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1811.1111
$ws.Range("I58").Value = 1425
$ws.Range("K58").Value = 1425
$ws.Range("M58").Value = -1222

# Row 94 (CRP) - hunk 40
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 2306
$ws.Range("I94").Value = 2306
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 2306
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -1855
$ws.Range("N94").Value = ""

# Row 103 (CRP) - hunk 41
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H103").Value = 125017940
$ws.Range("I103").Value = 200018300
$ws.Range("J103").Value = 17333.334
$ws.Range("K103").Value = 200018300
$ws.Range("L103").Value = 17333.334
$ws.Range("M103").Value = -200017128
$ws.Range("N103").Value = -19677.334

# Row 105 (CRP) - hunk 42
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 3862.7
$ws.Range("I105").Value = 4322.5713
$ws.Range("J105").Value = 3460.3125
$ws.Range("K105").Value = 4322.5713
$ws.Range("L105").Value = 3460.3125
$ws.Range("M105").Value = -2575.5713
$ws.Range("N105").Value = -6954.3125

# Row 109 (CRP) - hunk 43
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H109").Value = 46761.668
$ws.Range("J109").Value = 46761.668
$ws.Range("L109").Value = 46761.668
$ws.Range("N109").Value = -48841.668

# Row 116 (CRP) - hunk 44
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H116").Value = 53637.5
$ws.Range("J116").Value = 53637.5
$ws.Range("L116").Value = 53637.5
$ws.Range("N116").Value = -62815.5

# Row 128 (CRP) - hunk 45
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H128").Value = 99332
$ws.Range("J128").Value = 99332
$ws.Range("L128").Value = 99332
$ws.Range("N128").Value = -109292

# Row 132 (CRP) - hunk 46
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 744.5333000000001
$ws.Range("I132").Value = 552.1539
$ws.Range("K132").Value = 1656.4617
$ws.Range("M132").Value = 873.5382999999999

# Row 134 (CRP) - hunk 47
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2567.7693
$ws.Range("I134").Value = 2378.7778
$ws.Range("J134").Value = 2993
$ws.Range("K134").Value = 7136.3334
$ws.Range("L134").Value = 8979
$ws.Range("M134").Value = -4601.3334
$ws.Range("N134").Value = -14049

# Row 136 (CRP) - hunk 48
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1811.1111
$ws.Range("I136").Value = 1425
$ws.Range("K136").Value = 4275
$ws.Range("M136").Value = -1725

# Row 4 (CUL) - hunk 49
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 13480503
$ws.Range("I4").Value = 3766733.5
$ws.Range("K4").Value = 11300200.5
$ws.Range("M4").Value = -11300088.5

# Row 11 (CUL) - hunk 50
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 362.5484
$ws.Range("I11").Value = 351.3
$ws.Range("K11").Value = 1053.9
$ws.Range("M11").Value = -913.9000000000001

# Row 55 (CUL) - hunk 51
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 31251124
$ws.Range("I55").Value = 500
$ws.Range("J55").Value = 41668000
$ws.Range("K55").Value = 1500
$ws.Range("L55").Value = 125004000
$ws.Range("M55").Value = -1323
$ws.Range("N55").Value = -125004354

# Row 56 (CUL) - hunk 52
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 10254.714
$ws.Range("I56").Value = 10254.714
$ws.Range("K56").Value = 10254.714
$ws.Range("M56").Value = -9724.714

# Row 62 (CUL) - hunk 53
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 12666.333
$ws.Range("J62").Value = 14999.5
$ws.Range("L62").Value = 44998.5
$ws.Range("N62").Value = -46370.5

# Row 63 (CUL) - hunk 54
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").Value = ""

# Row 65 (CUL) - hunk 55
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H65").Value = 12666.333
$ws.Range("J65").Value = 14999.5
$ws.Range("L65").Value = 134995.5
$ws.Range("N65").Value = -141859.5

# Row 66 (CUL) - hunk 56
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").Value = ""

# Row 129 (CUL) - hunk 57
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 82840.56
$ws.Range("I129").Value = 500993.75
$ws.Range("J129").Value = 3192.3333
$ws.Range("K129").Value = 1502981.25
$ws.Range("L129").Value = 9576.999899999999
$ws.Range("M129").Value = -1497981.25
$ws.Range("N129").Value = -19576.9999

# Row 131 (CUL) - hunk 58
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1572.4231
$ws.Range("I131").Value = 1130
$ws.Range("J131").Value = 2175.7273
$ws.Range("K131").Value = 3390
$ws.Range("L131").Value = 6527.1819
$ws.Range("M131").Value = 1650
$ws.Range("N131").Value = -16607.1819

# Row 137 (CUL) - hunk 59
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 3229704.8
$ws.Range("I137").Value = 10002048
$ws.Range("J137").Value = 4779.381
$ws.Range("K137").Value = 30006144
$ws.Range("L137").Value = 14338.143
$ws.Range("M137").Value = -30001044
$ws.Range("N137").Value = -24538.143

# Row 2 (GSM) - hunk 60
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 534.2273
$ws.Range("I2").Value = 248.76923
$ws.Range("J2").Value = 946.55554
$ws.Range("K2").Value = 248.76923
$ws.Range("L2").Value = 946.55554
$ws.Range("M2").Value = -135.76923
$ws.Range("N2").Value = -1172.55554

# Row 44 (GSM) - hunk 61
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H44").Value = 25997.5
$ws.Range("J44").Value = 25997.5
$ws.Range("L44").Value = 25997.5
$ws.Range("N44").Value = -27189.5

# Row 99 (GSM) - hunk 62
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H99").Value = 13178.5
$ws.Range("J99").Value = 25259.5
$ws.Range("L99").Value = 25259.5
$ws.Range("N99").Value = -29751.5

# Row 102 (GSM) - hunk 63
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1411.9166
$ws.Range("I102").Value = 867.875
$ws.Range("K102").Value = 867.875
$ws.Range("M102").Value = 754.125

# Row 107 (GSM) - hunk 64
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 18909.46
$ws.Range("I107").Value = 27310.176
$ws.Range("K107").Value = 27310.176
$ws.Range("M107").Value = -25390.176

# Row 126 (GSM) - hunk 65
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2706.6
$ws.Range("I126").Value = 2151.1428
$ws.Range("J126").Value = 4002.6667
$ws.Range("K126").Value = 6453.428400000001
$ws.Range("L126").Value = 12008.0001
$ws.Range("M126").Value = -3983.428400000001
$ws.Range("N126").Value = -16948.0001

# Row 132 (GSM) - hunk 66
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1511.25
$ws.Range("I132").Value = 1512.6
$ws.Range("J132").Value = 1504.5
$ws.Range("K132").Value = 4537.799999999999
$ws.Range("L132").Value = 4513.5
$ws.Range("M132").Value = -2007.799999999999
$ws.Range("N132").Value = -9573.5

# Row 7 (LTW) - hunk 67
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 24305.154
$ws.Range("I7").Value = 36246
$ws.Range("J7").Value = 5199.8
$ws.Range("K7").Value = 36246
$ws.Range("L7").Value = 5199.8
$ws.Range("M7").Value = -36134
$ws.Range("N7").Value = -5423.8

# Row 40 (LTW) - hunk 68
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5403.647
$ws.Range("I40").Value = 4904.857
$ws.Range("K40").Value = 4904.857
$ws.Range("M40").Value = -4768.857

# Row 122 (LTW) - hunk 69
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3482.2144
$ws.Range("I122").Value = 3088.5557
$ws.Range("J122").Value = 4190.8
$ws.Range("K122").Value = 9265.667099999999
$ws.Range("L122").Value = 12572.4
$ws.Range("M122").Value = -6815.667099999999
$ws.Range("N122").Value = -17472.4

# Row 126 (LTW) - hunk 70
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 24305.154
$ws.Range("I126").Value = 36246
$ws.Range("J126").Value = 5199.8
$ws.Range("K126").Value = 108738
$ws.Range("L126").Value = 15599.4
$ws.Range("M126").Value = -106268
$ws.Range("N126").Value = -20539.4

# Row 136 (LTW) - hunk 71
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 3316.1875
$ws.Range("I136").Value = 2420.625
$ws.Range("K136").Value = 7261.875
$ws.Range("M136").Value = -4711.875

# Row 96 (WVR) - hunk 72
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1897.0834
$ws.Range("I96").Value = 1529.7778
$ws.Range("K96").Value = 1529.7778
$ws.Range("M96").Value = -156.7778000000001

# Row 122 (WVR) - hunk 73
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1295.909
$ws.Range("I122").Value = 1295.909
$ws.Range("K122").Value = 3887.727
$ws.Range("M122").Value = -1437.727

# Row 132 (WVR) - hunk 74
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 8735.682000000001
$ws.Range("I132").Value = 8709.25
$ws.Range("J132").Value = 9000
$ws.Range("K132").Value = 26127.75
$ws.Range("L132").Value = 27000
$ws.Range("M132").Value = -23597.75
$ws.Range("N132").Value = -32060
